# Estadisticos Segundo Parcial 26 Mayo
$wb = $excel.ActiveWorkbook

# Update "Estadisticos 2P" sheet with the new second-partial results.
$ws2P = $wb.Worksheets.Item("Estadisticos 2P")
$ws2P.Range("D2").Value = 0
$ws2P.Range("E2").Value = 1
$ws2P.Range("F2").Value = 34
$ws2P.Range("G2").Value = 97.14
$ws2P.Range("H2").Value = 7.4

# The recalculated final-average statistic also changes on "Estadisticos Final".
$wsFinal = $wb.Worksheets.Item("Estadisticos Final")
$wsFinal.Range("H2").Value = 7.3
